# Generate Report for Handoff
# The "fcee92e4-e213-4c85-bdcd-3f29e3f80fcf" file finished handoff, so its
# status moves from "In Translation" to "Ready for handoff" on every sheet
# that surfaces it, and the handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = "Ready for handoff"
$overview.Range("F4").Value = "Ready for handoff"
$overview.Range("G4").Value = "2016-10-10 09:30:12"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("H4").Value = "2016-10-10 09:30:00"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-10-10 09:30:12"
